# "added multiplayer to planning"
#
# Slide 2 ("Planning") has a content placeholder whose last bullet reads
# "Websockets". The edit appends a new bullet paragraph "Multiplayer"
# right after it, leaving every other paragraph/run untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the body/content placeholder shape on the slide (the one that is
# not the title) so the script does not depend on a brittle shape index.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.Name -ne "Titel 1") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$tr = $shape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
$lastPara = $tr.Paragraphs($paraCount, 1)

# Only add the new bullet if it is not already there (idempotent / safe).
if ($lastPara.Text -ne "Multiplayer") {
    $lastPara.InsertAfter("`rMultiplayer")
}
